$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.749.27"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "2.251.70"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "303.73"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "93.90"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").Value = "34.52"
$ws.Range("E10").Value = "  +5.26%  "
$ws.Range("D11").Value = "0.0785"
$ws.Range("E11").Value = "  -2.27%  "
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("D14").Value = "2.600.47"
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D15").Value = "14.25"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").Value = "2.257.65"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "41.647.96"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").Value = "12.23"
$ws.Range("E19").Value = "  -4.37%  "
$ws.Range("D20").Value = "0.0₃0895"
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("D21").Value = "5.92"
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("D22").Value = "67.67"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").Value = "235.59"
$ws.Range("E23").Value = "  -3.61%  "
$ws.Range("E24").Value = "  -2.39%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D27").Value = "23.50"
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("D28").Value = "35.85"
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("D30").Value = "9.38"
$ws.Range("E30").Value = "  -3.36%  "
$ws.Range("D31").Value = "158.91"
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  -3.65%  "
$ws.Range("E35").Value = "  -2.19%  "
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("D37").Value = "16.84"
$ws.Range("E37").Value = "  -2.73%  "
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("E40").Value = "  -3.02%  "
$ws.Range("D41").Value = "3.94"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").Value = "2.29"
$ws.Range("E42").Value = "  +2.66%  "
$ws.Range("D43").Value = "1.960.90"
$ws.Range("E43").Value = "  -2.76%  "
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("D45").Value = "18.55"
$ws.Range("E45").Value = "  -7.07%  "
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").Value = "9.80"
$ws.Range("E47").Value = "  -5.22%  "
$ws.Range("D48").Value = "52.50"
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("D51").Value = "90.31"
$ws.Range("E51").Value = "  -1.82%  "
